$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) / Volume(1h) (E) values scraped for this run.
# Values are forced to stay text (matching the original inlineStr cells)
# by temporarily switching the cell to a text number-format before the
# assignment (otherwise Excel auto-coerces strings like "317.80" or
# "0.00001074" into numbers and mangles/truncates them), then clearing
# the format back off so no stray style/quote-prefix is left behind.

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.029.32"
Set-TextValue "E2" "  -0.09%  "
Set-TextValue "D3" "1.868.97"
Set-TextValue "E3" "  -2.79%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "317.80"
Set-TextValue "E5" "  -2.26%  "
Set-TextValue "E6" "  -0.03%  "
Set-TextValue "D7" "0.5081"
Set-TextValue "E7" "  -1.51%  "
Set-TextValue "D8" "0.3917"
Set-TextValue "E8" "  -2.01%  "
Set-TextValue "D9" "0.08138"
Set-TextValue "E9" "  -3.89%  "
Set-TextValue "D10" "41.95"
Set-TextValue "E10" "  -2.37%  "
Set-TextValue "D11" "1.090"
Set-TextValue "E11" "  -2.88%  "
Set-TextValue "D12" "22.67"
Set-TextValue "E12" "  +7.12%  "
Set-TextValue "D13" "1.875.10"
Set-TextValue "E13" "  -2.33%  "
Set-TextValue "D14" "6.254"
Set-TextValue "E14" "  -1.16%  "
Set-TextValue "D15" "7.149"
Set-TextValue "E15" "  -2.59%  "
Set-TextValue "E16" "  +0.00%  "
Set-TextValue "D17" "91.58"
Set-TextValue "E17" "  -2.85%  "
Set-TextValue "D18" "0.00001074"
Set-TextValue "E18" "  -3.71%  "
Set-TextValue "D19" "0.06320"
Set-TextValue "E19" "  -6.49%  "
Set-TextValue "D20" "17.84"
Set-TextValue "E20" "  -0.85%  "
Set-TextValue "E21" "  +0.05%  "
Set-TextValue "D22" "30.020.90"
Set-TextValue "E22" "  -0.16%  "
Set-TextValue "D23" "5.779"
Set-TextValue "E23" "  -4.57%  "
Set-TextValue "D24" "11.05"
Set-TextValue "E24" "  -1.17%  "
Set-TextValue "D25" "2.202"
Set-TextValue "E25" "  -0.06%  "
Set-TextValue "D26" "2.082.65"
Set-TextValue "E26" "  -2.66%  "
Set-TextValue "D27" "160.49"
Set-TextValue "E27" "  +0.23%  "
Set-TextValue "D28" "20.82"
Set-TextValue "E28" "  -0.72%  "
Set-TextValue "D29" "2.221"
Set-TextValue "E29" "  -9.85%  "
Set-TextValue "D30" "126.29"
Set-TextValue "E30" "  -2.10%  "
Set-TextValue "D31" "0.1030"
Set-TextValue "E31" "  -2.56%  "
Set-TextValue "D32" "1.037"
Set-TextValue "E32" "  -3.88%  "
Set-TextValue "D33" "5.864"
Set-TextValue "E33" "  -3.39%  "
Set-TextValue "D34" "3.736"
Set-TextValue "E34" "  +2.07%  "
Set-TextValue "E35" "  -3.62%  "
Set-TextValue "D36" "5.174"
Set-TextValue "E36" "  -0.37%  "
Set-TextValue "D37" "0.06310"
Set-TextValue "E37" "  -4.45%  "
Set-TextValue "D38" "0.2131"
Set-TextValue "E38" "  -3.92%  "
Set-TextValue "D39" "1.170"
Set-TextValue "E39" "  -5.89%  "
Set-TextValue "D40" "8.459"
Set-TextValue "E40" "  -5.86%  "
Set-TextValue "D41" "0.6252"
Set-TextValue "E41" "  -4.28%  "
Set-TextValue "E42" "  -2.94%  "
Set-TextValue "D43" "11.22"
Set-TextValue "E43" "  -1.65%  "
Set-TextValue "E44" "  +0.00%  "
Set-TextValue "D45" "0.5854"
Set-TextValue "E45" "  -4.51%  "
Set-TextValue "D46" "12.77"
Set-TextValue "E46" "  -2.69%  "
Set-TextValue "D47" "3.625"
Set-TextValue "E47" "  -3.27%  "
Set-TextValue "D48" "1.982"
Set-TextValue "E48" "  -3.63%  "
Set-TextValue "D49" "121.78"
Set-TextValue "E49" "  -2.91%  "
Set-TextValue "D50" "1.198"
Set-TextValue "E50" "  -3.61%  "
Set-TextValue "D51" "1.146"
Set-TextValue "E51" "  +0.00%  "
